$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.280.92'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.648.16'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.59'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0638'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.95'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.877.32'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.661.53'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.550'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.50'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.281.95'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '196.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.34'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.05'
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +1.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.99'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.67'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E31").Value = '  +2.21%  '
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.25'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.61'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.31%  '
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.916'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.557'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.137.03'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.50'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.50'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.786.23'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.45'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0519'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.417'
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0972'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.89%  '
